{"js": "// Update the two-digit multiplication practice sheet with a new set of\n// problems. Each old \"AA\u00d7BB=\" expression in the table is replaced with its\n// corresponding new expression, in document order.\n\nconst replacements = [\n  [\"28\u00d783=\", \"98\u00d769=\"],\n  [\"96\u00d739=\", \"80\u00d783=\"],\n  [\"62\u00d736=\", \"40\u00d794=\"],\n  [\"22\u00d732=\", \"96\u00d795=\"],\n  [\"94\u00d741=\", \"29\u00d797=\"],\n  [\"52\u00d797=\", \"15\u00d714=\"],\n  [\"39\u00d784=\", \"33\u00d716=\"],\n  [\"68\u00d721=\", \"40\u00d760=\"],\n  [\"96\u00d714=\", \"28\u00d725=\"],\n  [\"53\u00d745=\", \"11\u00d788=\"],\n  [\"26\u00d712=\", \"80\u00d723=\"],\n  [\"49\u00d784=\", \"21\u00d727=\"],\n  [\"71\u00d796=\", \"45\u00d748=\"],\n  [\"18\u00d740=\", \"39\u00d715=\"],\n  [\"24\u00d765=\", \"18\u00d714=\"],\n  [\"17\u00d737=\", \"80\u00d739=\"],\n  [\"12\u00d722=\", \"73\u00d766=\"],\n  [\"69\u00d784=\", \"73\u00d763=\"],\n  [\"58\u00d760=\", \"92\u00d720=\"],\n  [\"38\u00d799=\", \"36\u00d798=\"],\n  [\"39\u00d734=\", \"82\u00d752=\"],\n  [\"29\u00d771=\", \"68\u00d723=\"],\n  [\"12\u00d754=\", \"49\u00d745=\"],\n  [\"47\u00d787=\", \"59\u00d732=\"],\n  [\"27\u00d757=\", \"61\u00d752=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication practice sheet with a new set of\n# problems. Each old \"AA\u00d7BB=\" expression in the table is replaced with its\n# corresponding new expression, in document order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"28\u00d783=\"; New = \"98\u00d769=\" },\n    @{ Old = \"96\u00d739=\"; New = \"80\u00d783=\" },\n    @{ Old = \"62\u00d736=\"; New = \"40\u00d794=\" },\n    @{ Old = \"22\u00d732=\"; New = \"96\u00d795=\" },\n    @{ Old = \"94\u00d741=\"; New = \"29\u00d797=\" },\n    @{ Old = \"52\u00d797=\"; New = \"15\u00d714=\" },\n    @{ Old = \"39\u00d784=\"; New = \"33\u00d716=\" },\n    @{ Old = \"68\u00d721=\"; New = \"40\u00d760=\" },\n    @{ Old = \"96\u00d714=\"; New = \"28\u00d725=\" },\n    @{ Old = \"53\u00d745=\"; New = \"11\u00d788=\" },\n    @{ Old = \"26\u00d712=\"; New = \"80\u00d723=\" },\n    @{ Old = \"49\u00d784=\"; New = \"21\u00d727=\" },\n    @{ Old = \"71\u00d796=\"; New = \"45\u00d748=\" },\n    @{ Old = \"18\u00d740=\"; New = \"39\u00d715=\" },\n    @{ Old = \"24\u00d765=\"; New = \"18\u00d714=\" },\n    @{ Old = \"17\u00d737=\"; New = \"80\u00d739=\" },\n    @{ Old = \"12\u00d722=\"; New = \"73\u00d766=\" },\n    @{ Old = \"69\u00d784=\"; New = \"73\u00d763=\" },\n    @{ Old = \"58\u00d760=\"; New = \"92\u00d720=\" },\n    @{ Old = \"38\u00d799=\"; New = \"36\u00d798=\" },\n    @{ Old = \"39\u00d734=\"; New = \"82\u00d752=\" },\n    @{ Old = \"29\u00d771=\"; New = \"68\u00d723=\" },\n    @{ Old = \"12\u00d754=\"; New = \"49\u00d745=\" },\n    @{ Old = \"47\u00d787=\"; New = \"59\u00d732=\" },\n    @{ Old = \"27\u00d757=\"; New = \"61\u00d752=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
